$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain Text so numeric-looking values (e.g. "1.00") are not
# silently converted to numbers by Excel's type inference on assignment.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.382.76"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "2.309.74"
$ws.Range("E3").Value = "  +2.42%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "311.19"
$ws.Range("E5").Value = "  +1.74%  "
$ws.Range("D6").Value = "103.31"
$ws.Range("E6").Value = "  +7.14%  "
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "0.530"
$ws.Range("E9").Value = "  +8.79%  "
$ws.Range("D10").Value = "36.48"
$ws.Range("E10").Value = "  +4.82%  "
$ws.Range("D11").Value = "0.0815"
$ws.Range("E11").Value = "  +4.07%  "
$ws.Range("D12").Value = "51.82"
$ws.Range("E12").Value = "  -0.69%  "
$ws.Range("E13").Value = "  -1.16%  "
$ws.Range("D14").Value = "7.05"
$ws.Range("E14").Value = "  +3.92%  "
$ws.Range("D15").Value = "2.665.70"
$ws.Range("E15").Value = "  +2.35%  "
$ws.Range("D16").Value = "15.08"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").Value = "2.310.40"
$ws.Range("E17").Value = "  +2.38%  "
$ws.Range("D18").Value = "0.810"
$ws.Range("E18").Value = "  +3.33%  "
$ws.Range("D19").Value = "43.281.46"
$ws.Range("E19").Value = "  +3.41%  "
$ws.Range("D20").Value = "12.35"
$ws.Range("E20").Value = "  +1.61%  "
$ws.Range("D21").Value = "0.0₃0932"
$ws.Range("E21").Value = "  +3.92%  "
$ws.Range("D22").Value = "6.20"
$ws.Range("E22").Value = "  +4.09%  "
$ws.Range("D23").Value = "68.14"
$ws.Range("E23").Value = "  +0.96%  "
$ws.Range("D24").Value = "241.93"
$ws.Range("E24").Value = "  +2.82%  "
$ws.Range("D25").Value = "2.03"
$ws.Range("E25").Value = "  +2.71%  "
$ws.Range("D26").Value = "2.62"
$ws.Range("E26").Value = "  +2.16%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("D28").Value = "24.84"
$ws.Range("E28").Value = "  +6.19%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "2.30"
$ws.Range("E29").Value = "  +8.40%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").Value = "36.82"
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "9.66"
$ws.Range("E31").Value = "  +1.75%  "
$ws.Range("D32").Value = "168.45"
$ws.Range("E32").Value = "  +2.23%  "
$ws.Range("D33").Value = "5.29"
$ws.Range("E33").Value = "  +1.86%  "
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("B35").Value = "Celestia"
$ws.Range("C35").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D35").Value = "18.02"
$ws.Range("E35").Value = "  +3.16%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D36").Value = "2.54"
$ws.Range("E36").Value = "  +5.93%  "
$ws.Range("D37").Value = "0.0743"
$ws.Range("E37").Value = "  +1.90%  "
$ws.Range("E38").Value = "  -2.28%  "
$ws.Range("E39").Value = "  +4.79%  "
$ws.Range("D40").Value = "0.106"
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "4.45"
$ws.Range("E41").Value = "  +8.83%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "0.116"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "2.44"
$ws.Range("E43").Value = "  +6.20%  "
$ws.Range("D44").Value = "0.0296"
$ws.Range("E44").Value = "  +6.29%  "
$ws.Range("D45").Value = "1.981.25"
$ws.Range("E45").Value = "  +1.99%  "
$ws.Range("D46").Value = "19.25"
$ws.Range("E46").Value = "  +1.28%  "
$ws.Range("D47").Value = "3.02"
$ws.Range("E47").Value = "  +4.96%  "
$ws.Range("E48").Value = "  +0.22%  "
$ws.Range("B49").Value = "MultiversX"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D49").Value = "55.90"
$ws.Range("E49").Value = "  +5.90%  "
$ws.Range("B50").Value = "HuobiToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D50").Value = "2.95"
$ws.Range("E50").Value = "  +2.10%  "
$ws.Range("D51").Value = "1.61"
$ws.Range("E51").Value = "  +11.19%  "
